$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "ID"
$ws.Range("D2").Value = "Name"
$ws.Range("E2").Value = "Value"
